$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new blood-pressure log entries (rows 9-24), mirroring the
# manual-entry formatting already used in the sheet (quote-prefixed
# dates/numbers so Excel keeps them as text instead of auto-converting).

# Row 9
$ws.Range("A9").Value = "'06/04/25"
$ws.Range("B9").Value = "123/76"
$ws.Range("C9").Value = "'58"
$ws.Range("D9").Value = "6:51PM"
[void]$ws.Range("A8").Copy()
[void]$ws.Range("A9").PasteSpecial(-4122)
[void]$ws.Range("B2").Copy()
[void]$ws.Range("B9:D9").PasteSpecial(-4122)

# Row 10
$ws.Range("A10").Value = "'06/05/25"
$ws.Range("B10").Value = "132/85"
$ws.Range("C10").Value = "'48"
$ws.Range("D10").Value = "5:00AM"
[void]$ws.Range("A8").Copy()
[void]$ws.Range("A10").PasteSpecial(-4122)
[void]$ws.Range("B2").Copy()
[void]$ws.Range("B10:C10").PasteSpecial(-4122)
[void]$ws.Range("D1").Copy()
[void]$ws.Range("D10").PasteSpecial(-4122)

# Row 11
$ws.Range("A11").Value = "'06/06/25"
$ws.Range("B11").Value = "132/82"
$ws.Range("C11").Value = "'50"
$ws.Range("D11").Value = "5:00AM"
[void]$ws.Range("A8").Copy()
[void]$ws.Range("A11").PasteSpecial(-4122)
[void]$ws.Range("B2").Copy()
[void]$ws.Range("B11:C11").PasteSpecial(-4122)
[void]$ws.Range("D1").Copy()
[void]$ws.Range("D11").PasteSpecial(-4122)

# Row 12
$ws.Range("A12").Value = "'06/06/25"
$ws.Range("B12").Value = "136/89"
$ws.Range("C12").Value = "'54"
$ws.Range("D12").Value = "5:00PM"
[void]$ws.Range("A8").Copy()
[void]$ws.Range("A12").PasteSpecial(-4122)
[void]$ws.Range("B2").Copy()
[void]$ws.Range("B12:C12").PasteSpecial(-4122)
[void]$ws.Range("D1").Copy()
[void]$ws.Range("D12").PasteSpecial(-4122)

# Row 13
$ws.Range("A13").Value = "'06/07/25"
$ws.Range("B13").Value = "138/89"
$ws.Range("C13").Value = "'49"
$ws.Range("D13").Value = "4:46AM"
[void]$ws.Range("A8").Copy()
[void]$ws.Range("A13").PasteSpecial(-4122)
[void]$ws.Range("B2").Copy()
[void]$ws.Range("B13:C13").PasteSpecial(-4122)
[void]$ws.Range("D1").Copy()
[void]$ws.Range("D13").PasteSpecial(-4122)

# Row 14
$ws.Range("A14").Value = "'06/07/25"
$ws.Range("B14").Value = "133/83"
$ws.Range("C14").Value = "'66"
$ws.Range("D14").Value = "5:40PM"
[void]$ws.Range("A8").Copy()
[void]$ws.Range("A14").PasteSpecial(-4122)
[void]$ws.Range("B2").Copy()
[void]$ws.Range("B14:D14").PasteSpecial(-4122)

# Row 15
$ws.Range("A15").Value = "'06/08/25"
$ws.Range("B15").Value = "136/85"
$ws.Range("C15").Value = "'54"
$ws.Range("D15").Value = "5:37AM"
[void]$ws.Range("A8").Copy()
[void]$ws.Range("A15").PasteSpecial(-4122)
[void]$ws.Range("B2").Copy()
[void]$ws.Range("B15:D15").PasteSpecial(-4122)

# Row 16
$ws.Range("A16").Value = "'06/09/25"
$ws.Range("B16").Value = "130/86"
$ws.Range("C16").Value = "'47"
$ws.Range("D16").Value = "4:50AM"
[void]$ws.Range("A8").Copy()
[void]$ws.Range("A16").PasteSpecial(-4122)
[void]$ws.Range("B2").Copy()
[void]$ws.Range("B16:D16").PasteSpecial(-4122)

# Row 17
$ws.Range("A17").Value = "'06/09/25"
$ws.Range("B17").Value = "133/87"
$ws.Range("C17").Value = "'50"
$ws.Range("D17").Value = "8:24PM"
[void]$ws.Range("A8").Copy()
[void]$ws.Range("A17").PasteSpecial(-4122)
[void]$ws.Range("B2").Copy()
[void]$ws.Range("B17:D17").PasteSpecial(-4122)

# Row 18
$ws.Range("A18").Value = "'06/10/25"
$ws.Range("B18").Value = "137/89"
$ws.Range("C18").Value = "'54"
$ws.Range("D18").Value = "5:30AM"
[void]$ws.Range("A8").Copy()
[void]$ws.Range("A18").PasteSpecial(-4122)
[void]$ws.Range("B2").Copy()
[void]$ws.Range("B18:D18").PasteSpecial(-4122)

# Row 19
$ws.Range("A19").Value = "'06/10/25"
$ws.Range("B19").Value = "126/76"
$ws.Range("C19").Value = "'59"
$ws.Range("D19").Value = "5:55PM"
[void]$ws.Range("A8").Copy()
[void]$ws.Range("A19").PasteSpecial(-4122)
[void]$ws.Range("B2").Copy()
[void]$ws.Range("B19:D19").PasteSpecial(-4122)

# Row 20
$ws.Range("A20").Value = "'06/11/25"
$ws.Range("B20").Value = "139/86"
$ws.Range("C20").Value = "'49"
$ws.Range("D20").Value = "5:12AM"
[void]$ws.Range("A8").Copy()
[void]$ws.Range("A20").PasteSpecial(-4122)
[void]$ws.Range("B2").Copy()
[void]$ws.Range("B20:D20").PasteSpecial(-4122)

# Row 21
$ws.Range("A21").Value = "'06/11/25"
$ws.Range("B21").Value = "134/88"
$ws.Range("D21").Value = "6:23PM"
[void]$ws.Range("A8").Copy()
[void]$ws.Range("A21").PasteSpecial(-4122)
[void]$ws.Range("B2").Copy()
[void]$ws.Range("B21:D21").PasteSpecial(-4122)

# Row 22
$ws.Range("A22").Value = "'06/12/25"
$ws.Range("B22").Value = "137/85"
$ws.Range("C22").Value = "'48"
$ws.Range("D22").Value = "5:00AM"
[void]$ws.Range("A8").Copy()
[void]$ws.Range("A22").PasteSpecial(-4122)
[void]$ws.Range("B2").Copy()
[void]$ws.Range("B22:D22").PasteSpecial(-4122)

# Row 23
$ws.Range("A23").Value = "'06/12/25"
$ws.Range("B23").Value = "133/83"
$ws.Range("C23").Value = "'55"
$ws.Range("D23").Value = "6:33PM"
[void]$ws.Range("A8").Copy()
[void]$ws.Range("A23").PasteSpecial(-4122)
[void]$ws.Range("B2").Copy()
[void]$ws.Range("B23:D23").PasteSpecial(-4122)

# Row 24
$ws.Range("D24").Value = "GGNN,.."
[void]$ws.Range("D1").Copy()
[void]$ws.Range("D24").PasteSpecial(-4122)

# Leave the selection where the author last left it before saving.
[void]$ws.Range("M23").Select()
